$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume(1h) columns hold plain display text (e.g. "71.098.93",
# "  -2.27%  "), not numbers. Force each touched cell to text storage
# right before writing so Excel does not reinterpret the string as a
# number and silently drop formatting like trailing zeros.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.980.79"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.945.64"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.82%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.32"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.98"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.940.49"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.77%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.687"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -5.69%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.739"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.96%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -6.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.38"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +14.45%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.61"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.577.64"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.949.53"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.89"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.56"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.65%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.16"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -4.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.971.51"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "420.18"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -7.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "97.63"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -7.33%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.63%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +5.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.54"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.27"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.84"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +16.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.71"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.64%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.47"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.70%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +16.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.89"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +19.92%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.34"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "678.61"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "65.71"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.440"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0816"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -6.13%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.92%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.71%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0482"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.77%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.17"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.62%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -5.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.66"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.24%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.00"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "144.06"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.44%  "
